$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, which pushes the existing rows
# 119-155 down to 120-156 (dates/prices shift down by one row, as seen
# in the diff).
$ws.Rows("119:119").Insert()

# Populate the newly inserted row 119 with the new weekly data entry.
$ws.Cells.Item(119, 1).Value = 8
$ws.Cells.Item(119, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44468
$ws.Cells.Item(119, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(119, 5).Value = 4
$ws.Cells.Item(119, 6).Value = 100112012
$ws.Cells.Item(119, 7).Value = "Espinaca"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 2800
$ws.Cells.Item(119, 11).Value = 450
$ws.Cells.Item(119, 12).Value = 500
$ws.Cells.Item(119, 13).Value = 475
$ws.Cells.Item(119, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(119, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(119, 16).Value = 950
$ws.Cells.Item(119, 17).Value = 0.5
$ws.Cells.Item(119, 18).Value = "Hortaliza"
